$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor values updated
$ws.Range("B3").Value = 0.99539734596513
$ws.Range("C3").Value = 0.9952256443579341
$ws.Range("D3").Value = 0.9893748654215145

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9897778293470486
$ws.Range("C4").Value = 0.9898259249882889
$ws.Range("D4").Value = 0.9670764166799234

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9980523453027329
$ws.Range("C5").Value = 0.9980006692028729
$ws.Range("D5").Value = 0.9973938159072008
